# Auto-update draw results: append the 2025-11-27 Pick 4 draw as a new
# row at the bottom of the results table (row 72).
#
# The sheet stores every column (including the date-like "Date"/"Phase"
# columns) as literal TEXT, not as real dates/numbers. Setting
# Range.Value directly on a string like "2025-11-27" makes Excel infer a
# date/number and silently reformat the cell (and would also stamp a new
# "quote prefix" style on it). To avoid that, we write each value as a
# text-formula (="...") first - which always yields a string result -
# and then flatten it back down to a plain value with
# Copy / PasteSpecial(xlPasteValues), exactly like using Excel's
# "Paste Values" after typing a formula. That keeps the cells as plain
# text with no extra formula and no style changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRowNum = 72
$rowValues = @(
    "2025-11-27",
    "Pick 4",
    "251127",
    "3-0-8-6",
    "2025-11-27T21:37:50.498+04:00"
)

for ($i = 0; $i -lt $rowValues.Count; $i++) {
    $colNum = $i + 1
    $ws.Cells.Item($newRowNum, $colNum).Formula = '="' + $rowValues[$i] + '"'
}

$newRowRange = $ws.Range($ws.Cells.Item($newRowNum, 1), $ws.Cells.Item($newRowNum, $rowValues.Count))
$newRowRange.Copy()
$newRowRange.PasteSpecial(-4163)   # xlPasteValues - collapse the helper formulas into literal text
$excel.CutCopyMode = $false
